$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row -> new value for column F (dSF)
$updates = @{
    2  = -1
    6  = -3
    9  = 2
    10 = -1
    18 = 1
    19 = 1
    22 = 4
    29 = 6
    30 = -1
    31 = -1
    35 = -1
    63 = -2
    68 = -5
    70 = 2
    75 = -1
    79 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
